$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: "15. Dynamic Programming" / "Longest Palindrome Substring" (Neutral style) / hyperlink-text url (no live hyperlink)
$ws.Range("D15").Value = "https://leetcode.com/problems/longest-palindromic-substring/"
$ws.Range("B15").Value = "Longest Palindrome Substring"
$ws.Range("B15").Style = "Neutral"
$ws.Range("A15").Value = "15. Dynamic Programming"

# Row 16: "16. Arrays" / "Running Sum of 1d Array" (Good style) / url text
$ws.Range("D16").Value = "https://leetcode.com/problems/running-sum-of-1d-array/"
$ws.Range("B16").Value = "Running Sum of 1d Array"
$ws.Range("B16").Style = "Good"
$ws.Range("A16").Value = "16. Arrays"

# Move the active selection like the author left it
$ws.Range("C24").Select() | Out-Null
